$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.714.75'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '3.745.50'
$ws.Range("E3").Value = '  +6.34%  '
$ws.Range("E4").Value = '  -0.12%  '
$c = $ws.Range("D5")
$c.Value = "'613.15"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +3.97%  '
$c = $ws.Range("D6")
$c.Value = "'177.92"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -4.04%  '
$ws.Range("D7").Value = '3.741.33'
$ws.Range("E7").Value = '  +6.26%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("E10").Value = '  +5.37%  '
$c = $ws.Range("D11")
$c.Value = "'6.35"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -3.40%  '
$c = $ws.Range("D12")
$c.Value = "'0.499"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.04%  '
$c = $ws.Range("D13")
$c.Value = "'40.92"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +5.97%  '
$ws.Range("E14").Value = '  +1.82%  '
$ws.Range("D15").Value = '4.365.60'
$ws.Range("E15").Value = '  +6.22%  '
$ws.Range("D16").Value = '3.740.52'
$ws.Range("E16").Value = '  +6.55%  '
$ws.Range("D17").Value = '69.786.16'
$ws.Range("E17").Value = '  +0.14%  '
$c = $ws.Range("D18")
$c.Value = "'0.124"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.44%  '
$c = $ws.Range("D19")
$c.Value = "'7.61"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.49%  '
$c = $ws.Range("D20")
$c.Value = "'515.33"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.40%  '
$c = $ws.Range("D21")
$c.Value = "'16.71"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.76%  '
$c = $ws.Range("D22")
$c.Value = "'9.59"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +7.64%  '
$ws.Range("E23").Value = '  -0.55%  '
$c = $ws.Range("D24")
$c.Value = "'88.13"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("E25").Value = '  +5.41%  '
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("E28").Value = '  +0.00%  '
$c = $ws.Range("D29")
$c.Value = "'0.0000127"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +18.43%  '
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("E31").Value = '  +4.15%  '
$ws.Range("E32").Value = '  -3.68%  '
$c = $ws.Range("D33")
$c.Value = "'31.36"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("E34").Value = '  -0.85%  '
$ws.Range("E35").Value = '  -0.19%  '
$c = $ws.Range("D36")
$c.Value = "'6.23"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.49%  '
$c = $ws.Range("D37")
$c.Value = "'1.03"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.35%  '
$ws.Range("E38").Value = '  +2.78%  '
$ws.Range("E39").Value = '  +3.54%  '
$c = $ws.Range("D40")
$c.Value = "'0.134"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +4.65%  '
$c = $ws.Range("D41")
$c.Value = "'51.23"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.08%  '
$c = $ws.Range("D42")
$c.Value = "'44.50"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -6.60%  '
$c = $ws.Range("D43")
$c.Value = "'8.84"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.43%  '
$c = $ws.Range("D44")
$c.Value = "'424.89"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.75%  '
$ws.Range("D45").Value = '3.093.50'
$ws.Range("E45").Value = '  +3.35%  '
$c = $ws.Range("D46")
$c.Value = "'2.73"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -3.26%  '
$c = $ws.Range("D47")
$c.Value = "'0.0365"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.33%  '
$c = $ws.Range("D48")
$c.Value = "'27.84"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("E49").Value = '  +3.80%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D50")
$c.Value = "'135.76"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("B51").Value = 'USDe'
$ws.Range("C51").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range("D51")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
